# Update mode: append newly approved vendor / new-hire VC records.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 - new vendor number
$ws.Cells.Item(43, 1).Value = "'1006914830"
$ws.Cells.Item(43, 2).Value = "Za Construction, LLC"

# Row 44 - new vendor number
$ws.Cells.Item(44, 1).Value = "'67902"
$ws.Cells.Item(44, 2).Value = "Koren Development Company"
